$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 500
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").Value = $null
$ws.Range("H121").Value = 5870.5835
$ws.Range("I121").Value = 750
$ws.Range("J121").Value = 6336.091
$ws.Range("K121").Value = 2250
$ws.Range("L121").Value = 19008.273
$ws.Range("M121").Value = -503
$ws.Range("N121").Value = -22502.273
$ws.Range("H129").Value = 1169.7222
$ws.Range("J129").Value = 1435.1852
$ws.Range("L129").Value = 4305.5556
$ws.Range("N129").Value = -14305.5556
$ws.Range("H132").Value = 3589.1035
$ws.Range("I132").Value = 3967.4092
$ws.Range("K132").Value = 11902.2276
$ws.Range("M132").Value = -9372.2276
$ws.Range("H135").Value = 45456172
$ws.Range("I135").Value = 1828.5714
$ws.Range("J135").Value = 125001260
$ws.Range("K135").Value = 16457.1426
$ws.Range("L135").Value = 1125011340
$ws.Range("M135").Value = -13922.1426
$ws.Range("N135").Value = -1125016410
$ws.Range("H137").Value = 1753.0625
$ws.Range("I137").Value = 1542
$ws.Range("J137").Value = 2667.6667
$ws.Range("K137").Value = 4626
$ws.Range("L137").Value = 8003.000100000001
$ws.Range("M137").Value = -2076
$ws.Range("N137").Value = -13103.0001
$ws.Range("H141").Value = 2106.6296
$ws.Range("I141").Value = 1408.5238
$ws.Range("K141").Value = 4225.5714
$ws.Range("M141").Value = 954.4286000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2609.1924
$ws.Range("I122").Value = 2226.75
$ws.Range("K122").Value = 6680.25
$ws.Range("M122").Value = -4230.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2745.9607
$ws.Range("I31").Value = 2006.8695
$ws.Range("J31").Value = 3353.0715
$ws.Range("K31").Value = 2006.8695
$ws.Range("L31").Value = 3353.0715
$ws.Range("M31").Value = -1711.8695
$ws.Range("N31").Value = -3943.0715
$ws.Range("H34").Value = 2745.9607
$ws.Range("I34").Value = 2006.8695
$ws.Range("J34").Value = 3353.0715
$ws.Range("K34").Value = 2006.8695
$ws.Range("L34").Value = 3353.0715
$ws.Range("M34").Value = -1804.8695
$ws.Range("N34").Value = -3757.0715
$ws.Range("H41").Value = 1059
$ws.Range("I41").Value = 1059
$ws.Range("K41").Value = 1059
$ws.Range("M41").Value = -631
$ws.Range("H50").Value = 20000
$ws.Range("J50").Value = 20000
$ws.Range("L50").Value = 20000
$ws.Range("N50").Value = -21250
$ws.Range("H51").Value = 20000
$ws.Range("J51").Value = 20000
$ws.Range("L51").Value = 20000
$ws.Range("N51").Value = -21472
$ws.Range("H60").Value = 14953.333
$ws.Range("J60").Value = 18657.143
$ws.Range("L60").Value = 18657.143
$ws.Range("N60").Value = -19679.143
$ws.Range("H61").Value = 20000
$ws.Range("J61").Value = 20000
$ws.Range("L61").Value = 20000
$ws.Range("N61").Value = -20696
$ws.Range("H68").Value = 49400
$ws.Range("J68").Value = 49400
$ws.Range("L68").Value = 49400
$ws.Range("N68").Value = -50898
$ws.Range("H71").Value = 49400
$ws.Range("J71").Value = 49400
$ws.Range("L71").Value = 148200
$ws.Range("N71").Value = -155688
$ws.Range("H74").Value = 39953.25
$ws.Range("J74").Value = 39953.25
$ws.Range("L74").Value = 39953.25
$ws.Range("N74").Value = -41701.25
$ws.Range("H77").Value = 39953.25
$ws.Range("J77").Value = 39953.25
$ws.Range("L77").Value = 119859.75
$ws.Range("N77").Value = -128595.75
$ws.Range("H99").Value = 26319358
$ws.Range("I99").Value = 2966.4167
$ws.Range("K99").Value = 2966.4167
$ws.Range("M99").Value = -1468.4167
$ws.Range("H105").Value = 13889881
$ws.Range("I105").Value = 25000522
$ws.Range("J105").Value = 1580.5
$ws.Range("K105").Value = 25000522
$ws.Range("L105").Value = 1580.5
$ws.Range("M105").Value = -24998775
$ws.Range("N105").Value = -5074.5
$ws.Range("H126").Value = 26319358
$ws.Range("I126").Value = 2966.4167
$ws.Range("K126").Value = 8899.250100000001
$ws.Range("M126").Value = -6429.250100000001
$ws.Range("H132").Value = 2702.9
$ws.Range("I132").Value = 2074.3044
$ws.Range("J132").Value = 4768.2856
$ws.Range("K132").Value = 6222.9132
$ws.Range("L132").Value = 14304.8568
$ws.Range("M132").Value = -3692.9132
$ws.Range("N132").Value = -19364.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 647.5
$ws.Range("I14").Value = 647.5
$ws.Range("K14").Value = 1942.5
$ws.Range("M14").Value = -1769.5
$ws.Range("H80").Value = 1215.7368
$ws.Range("I80").Value = 1949.5
$ws.Range("J80").Value = 1129.4117
$ws.Range("K80").Value = 5848.5
$ws.Range("L80").Value = 3388.2351
$ws.Range("M80").Value = -4912.5
$ws.Range("N80").Value = -5260.2351
$ws.Range("H83").Value = 1215.7368
$ws.Range("I83").Value = 1949.5
$ws.Range("J83").Value = 1129.4117
$ws.Range("K83").Value = 17545.5
$ws.Range("L83").Value = 10164.7053
$ws.Range("M83").Value = -12865.5
$ws.Range("N83").Value = -19524.7053
$ws.Range("H131").Value = 744.3
$ws.Range("J131").Value = 780.67035
$ws.Range("L131").Value = 2342.01105
$ws.Range("N131").Value = -12422.01105
$ws.Range("H136").Value = 2984.4285
$ws.Range("I136").Value = 1597.5
$ws.Range("J136").Value = 4833.6665
$ws.Range("K136").Value = 4792.5
$ws.Range("L136").Value = 14500.9995
$ws.Range("M136").Value = 307.5
$ws.Range("N136").Value = -24700.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1332.9048
$ws.Range("I97").Value = 1388.125
$ws.Range("J97").Value = 1156.2
$ws.Range("K97").Value = 1388.125
$ws.Range("L97").Value = 1156.2
$ws.Range("M97").Value = -892.125
$ws.Range("N97").Value = -2148.2
$ws.Range("H132").Value = 17846.322
$ws.Range("I132").Value = 1251.091
$ws.Range("K132").Value = 3753.273
$ws.Range("M132").Value = -1223.273

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 983244.9
$ws.Range("I122").Value = 1403021.2
$ws.Range("K122").Value = 4209063.6
$ws.Range("M122").Value = -4206613.6
$ws.Range("H136").Value = 1875
$ws.Range("I136").Value = 1875
$ws.Range("K136").Value = 5625
$ws.Range("M136").Value = -3075

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 994.0769
$ws.Range("I122").Value = 967.26086
$ws.Range("K122").Value = 2901.78258
$ws.Range("M122").Value = -451.7825800000001
